# The two "Caro White Cream" sale lines (row 2: qty 2 / amt 3600, and the
# now-removed row 3: qty 5 / amt 9000, both dated 2023-10-29) were merged
# into a single corrected line on row 2: qty 3 / amt 5400, dated 2023-10-31.
# Row 3 is then deleted entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DATE/QTY/AMT on row 2. These columns store plain text (shared
# strings) in the workbook, so force Text format before writing the value
# to stop Excel from auto-converting "2023-10-31" to a date serial or "3"
# / "5400" to numbers. Reset the style back to Normal afterwards so no
# stray number-format style is left behind on the cell.
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "2023-10-31"
$ws.Cells.Item(2, 2).Style = "Normal"

$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "3"
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(2, 6).NumberFormat = "@"
$ws.Cells.Item(2, 6).Value = "5400"
$ws.Cells.Item(2, 6).Style = "Normal"

# Remove the now-redundant third row (shifts nothing below it up, it was
# the last row) and shrink the used range accordingly.
$ws.Range("A3:F3").Delete()
